$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 447886.4
$ws.Range("I76").Value = 743142.7
$ws.Range("J76").Value = 5002
$ws.Range("K76").Value = 743142.7
$ws.Range("L76").Value = 5002
$ws.Range("M76").Value = -742827.7
$ws.Range("N76").Value = -5632

$ws.Range("H79").Value = 447886.4
$ws.Range("I79").Value = 743142.7
$ws.Range("J79").Value = 5002
$ws.Range("K79").Value = 743142.7
$ws.Range("L79").Value = 5002
$ws.Range("M79").Value = -742050.7
$ws.Range("N79").Value = -7186

$ws.Range("H106").Value = 976.25
$ws.Range("I106").Value = 968.3333
$ws.Range("J106").Value = 1000
$ws.Range("K106").Value = 968.3333
$ws.Range("L106").Value = 1000
$ws.Range("M106").Value = -337.3333
$ws.Range("N106").Value = -2262

$ws.Range("H118").Value = 4857.4243
$ws.Range("I118").Value = 758.5333000000001
$ws.Range("J118").Value = 8273.166999999999
$ws.Range("K118").Value = 2275.5999
$ws.Range("L118").Value = 24819.501
$ws.Range("M118").Value = -618.5999000000002
$ws.Range("N118").Value = -28133.501

$ws.Range("H138").Value = 4714.56
$ws.Range("I138").Value = 1623.2
$ws.Range("J138").Value = 5487.4
$ws.Range("K138").Value = 4869.6
$ws.Range("L138").Value = 16462.2
$ws.Range("M138").Value = 270.3999999999996
$ws.Range("N138").Value = -26742.2

$ws.Range("H141").Value = 3003.4243
$ws.Range("I141").Value = 1773.5769
$ws.Range("J141").Value = 7571.4287
$ws.Range("K141").Value = 5320.7307
$ws.Range("L141").Value = 22714.2861
$ws.Range("M141").Value = -140.7307000000001
$ws.Range("N141").Value = -33074.2861

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 29076.79
$ws.Range("I32").Value = 15941.6
$ws.Range("J32").Value = 42211.98
$ws.Range("K32").Value = 15941.6
$ws.Range("L32").Value = 42211.98
$ws.Range("M32").Value = -15654.6
$ws.Range("N32").Value = -42785.98

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 287079.1
$ws.Range("I134").Value = 323574.12
$ws.Range("J134").Value = 4242.5
$ws.Range("K134").Value = 970722.36
$ws.Range("L134").Value = 12727.5
$ws.Range("M134").Value = -968187.36
$ws.Range("N134").Value = -17797.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2591.574
$ws.Range("I31").Value = 1390.1923
$ws.Range("J31").Value = 3707.1428
$ws.Range("K31").Value = 1390.1923
$ws.Range("L31").Value = 3707.1428
$ws.Range("M31").Value = -1095.1923
$ws.Range("N31").Value = -4297.1428

$ws.Range("H34").Value = 2591.574
$ws.Range("I34").Value = 1390.1923
$ws.Range("J34").Value = 3707.1428
$ws.Range("K34").Value = 1390.1923
$ws.Range("L34").Value = 3707.1428
$ws.Range("M34").Value = -1188.1923
$ws.Range("N34").Value = -4111.1428

$ws.Range("H107").Value = 2062.4
$ws.Range("I107").Value = 3699.6667
$ws.Range("J107").Value = 1360.7142
$ws.Range("K107").Value = 3699.6667
$ws.Range("L107").Value = 1360.7142
$ws.Range("M107").Value = -1779.6667
$ws.Range("N107").Value = -5200.7142

$ws.Range("H132").Value = 2993.3572
$ws.Range("I132").Value = 2513.5
$ws.Range("J132").Value = 3353.25
$ws.Range("K132").Value = 7540.5
$ws.Range("L132").Value = 10059.75
$ws.Range("M132").Value = -5010.5
$ws.Range("N132").Value = -15119.75

$ws.Range("H141").Value = 52129.41
$ws.Range("I141").Value = 31765.334
$ws.Range("J141").Value = 56493.145
$ws.Range("K141").Value = 31765.334
$ws.Range("L141").Value = 56493.145
$ws.Range("M141").Value = -26585.334
$ws.Range("N141").Value = -66853.14499999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 667.0635
$ws.Range("I5").Value = 506.69388
$ws.Range("J5").Value = 1228.3572
$ws.Range("K5").Value = 1520.08164
$ws.Range("L5").Value = 3685.0716
$ws.Range("M5").Value = -1408.08164
$ws.Range("N5").Value = -3909.0716

$ws.Range("H32").Value = 2251
$ws.Range("J32").Value = 2500
$ws.Range("L32").Value = 7500
$ws.Range("N32").Value = -8066

$ws.Range("H46").Value = 2715.6843
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 2715.6843
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 8147.0529
$ws.Range("M46").ClearContents()
$ws.Range("N46").Value = -8329.052899999999

$ws.Range("H122").Value = 1024.5518
$ws.Range("I122").Value = 437.33334
$ws.Range("K122").Value = 3936.00006
$ws.Range("M122").Value = -1486.00006

$ws.Range("H131").Value = 1516068.2
$ws.Range("I131").Value = 4167191.5
$ws.Range("J131").Value = 1140.7858
$ws.Range("K131").Value = 12501574.5
$ws.Range("L131").Value = 3422.3574
$ws.Range("M131").Value = -12496534.5
$ws.Range("N131").Value = -13502.3574

$ws.Range("H135").Value = 667.0635
$ws.Range("I135").Value = 506.69388
$ws.Range("J135").Value = 1228.3572
$ws.Range("K135").Value = 4560.24492
$ws.Range("L135").Value = 11055.2148
$ws.Range("M135").Value = -2025.24492
$ws.Range("N135").Value = -16125.2148

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1368.8235
$ws.Range("I102").Value = 1411.1428
$ws.Range("J102").Value = 1171.3334
$ws.Range("K102").Value = 1411.1428
$ws.Range("L102").Value = 1171.3334
$ws.Range("M102").Value = 210.8571999999999
$ws.Range("N102").Value = -4415.3334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1961481.2
$ws.Range("J22").Value = 831.25
$ws.Range("L22").Value = 831.25
$ws.Range("N22").Value = -1421.25

$ws.Range("H27").Value = 1961481.2
$ws.Range("J27").Value = 831.25
$ws.Range("L27").Value = 831.25
$ws.Range("N27").Value = -1045.25

$ws.Range("H109").Value = 11136
$ws.Range("I109").Value = 9259
$ws.Range("J109").Value = 11761.667
$ws.Range("K109").Value = 9259
$ws.Range("L109").Value = 11761.667
$ws.Range("M109").Value = -7872
$ws.Range("N109").Value = -14535.667

$ws.Range("H122").Value = 4658.543
$ws.Range("I122").Value = 6042.45
$ws.Range("J122").Value = 2813.3333
$ws.Range("K122").Value = 18127.35
$ws.Range("L122").Value = 8439.999899999999
$ws.Range("M122").Value = -15677.35
$ws.Range("N122").Value = -13339.9999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H125").Value = 18992.75
$ws.Range("J125").Value = 18992.75
$ws.Range("L125").Value = 18992.75
$ws.Range("N125").Value = -28832.75
